$d = $word.ActiveDocument

# Replace WIFI_SSID value "NETGEAR64" with "rishabh"
$d.Content.Find.Execute("NETGEAR64", $true, $false, $false, $false, $false,
                         $true, 1, $false, "rishabh", 2)

# Replace WIFI_PASSWORD value "*JAVATPOINT#" with "rishabh"
$d.Content.Find.Execute("*JAVATPOINT#", $true, $false, $false, $false, $false,
                         $true, 1, $false, "rishabh", 2)
